$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D11 value changes from "unimportant" to "stuff"
$ws.Range("D11").Value = "stuff"

# D10 and D12 cells are cleared (removed entirely)
$ws.Range("D10").ClearContents()
$ws.Range("D12").ClearContents()

# Update the selection to D12
$ws.Range("D12").Select()
